$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each assignment uses a leading apostrophe to force text interpretation
# (matches the source file, where every cell is stored as inline/shared text,
# e.g. "206.29" must not become the Number 206.29). The Style reset immediately
# after clears the quotePrefix formatting flag that Value-assignment leaves behind,
# so the cell keeps its original (default) style index.
$ws.Range("D2").Value = "'26.907.27"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.05%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.549.78"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.16%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.33%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'206.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.13%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.488"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.74%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.34%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'22.04"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +2.27%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -0.12%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.0588"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +1.04%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0854"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.41%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.770.66"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.13%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.529.63"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.38%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  +0.99%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  +0.85%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'26.913.53"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.01%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'61.68"
$ws.Range("D17").Style = "Normal"
$ws.Range("B18").Value = "'BitcoinCash"
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").Value = "'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").Value = "'217.33"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.64%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("B19").Value = "'ShibaInu"
$ws.Range("B19").Style = "Normal"
$ws.Range("C19").Value = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("C19").Style = "Normal"
$ws.Range("D19").Value = "'0.0₃0703"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +2.66%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  +0.16%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E22").Value = "'  +0.67%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +0.63%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  -1.04%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'153.80"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.54%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'6.65"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.16%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  +0.59%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +0.92%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -0.24%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.0469"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +2.12%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -0.42%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -0.07%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'3.11"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +4.91%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.408.01"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +2.32%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.60"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +2.71%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.967"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.35%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +0.04%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +1.03%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.528"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +0.96%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.806"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.21%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -0.34%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  +3.49%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'2.30"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +2.06%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  +0.85%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'64.49"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +1.27%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  +0.45%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'1.684.96"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.12%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'87.43"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +1.36%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  +1.57%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  +3.55%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0960"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.61%  "
$ws.Range("E51").Style = "Normal"
